$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Provider Suvery Data")

$ws.Range("AX1").Value = '1 - Please describe your facility''s emergency communication infrastructure.'
$ws.Range("AY1").Value = 'What kind of communication connectivity is available at your site? - 1 - Cell Service'
$ws.Range("AZ1").Value = 'What kind of communication connectivity is available at your site? - 1 - Wi-Fi/Internet'
$ws.Range("BA1").Value = 'What kind of communication connectivity is available at your site? - 1 - Land Line'
$ws.Range("BB1").Value = 'Camp access: (check all that apply) - 1 - Camp access: (check all that apply) - Selected Choice'
$ws.Range("BC1").Value = 'Camp access: (check all that apply) - 1 - Other - Text'
$ws.Range("BD1").Value = '1 - Location''s acreage (if unknown please state unknown)'
$ws.Range("BE1").Value = '1 - What natural resources exist for study? Forest, Savanna, Stream, Pond? Special features such as salmon spawning location or wildfire evidence or unique geological features for study? Please describe.'
$ws.Range("BF1").Value = 'At (insert facility name), what food service aspects are provided? - 1 - At (insert facility name), what food service aspects are provided? - Selected Choice'
$ws.Range("BG1").Value = 'At (insert facility name), what food service aspects are provided? - 1 - Other - Text'
$ws.Range("BH1").Value = 'If food is provided, can menus accommodate special dietary needs? (check all that apply) - 1 - If food is provided, can menus accommodate special dietary needs? (check all that apply) - Selected Choice'
$ws.Range("BI1").Value = 'If food is provided, can menus accommodate special dietary needs? (check all that apply) - 1 - Other - Text'
$ws.Range("BJ1").Value = '1 - Maximum inside dining capacity? (if unknown please state unknown)'
$ws.Range("BK1").Value = '1 - Maximum outside dining capacity?  (if unknown please state unknown)'
$ws.Range("BL1").Value = 'Are there covered or indoor meeting areas? - 1 - Are there covered or indoor meeting areas? - Selected Choice'
$ws.Range("BM1").Value = 'Are there covered or indoor meeting areas? - 1 - If Yes, How many and what capacity? - Text'
$ws.Range("BN1").Value = '1 - Maximum overnight indoor sleeping capacity (please describe accommodations)'
$ws.Range("BO1").Value = '1 - Maximum overnight outside sleeping capacity (please describe accomodations)'
$ws.Range("BP1").Value = '1 - Does your facility offer universally accessible sleeping facilities to accommodate students with limited mobility? If so, please describe.'
$ws.Range("BQ1").Value = '1 - What are the biggest unmet facility and material needs for your outdoor school?'
$ws.Range("BR1").Value = '1 - Is there any additional information that you would like to provide regarding facilities and 
materials?'

$ws.Rows.Item(1).EntireRow.AutoFit()
